$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Tipo" column (D) to hold MAE values.
$ws.Columns.Item(4).Insert()

# Copy the header formatting (bold, border, centered) from the neighboring
# "R2" header cell onto the new "MAE" header cell.
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)

# Header text for the new column
$ws.Cells.Item(1, 4).Value = "MAE"

# MAE values for each row
$ws.Cells.Item(2, 4).Value = 0.9239883288152645
$ws.Cells.Item(3, 4).Value = 0.8429226117853587
$ws.Cells.Item(4, 4).Value = 0.886280647950258
$ws.Cells.Item(5, 4).Value = 0.8956978851546247
